$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Measured Watts [W]" / "Measured Voltage [mV]" column headers ---
$ws.Range("B1").Value = "Measured Voltage [mV]"
$ws.Range("C1").Value = "Measured Watts [W]"

# --- Swap the measured data between column B and column C (rows 2-8) ---
for ($r = 2; $r -le 8; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# --- Re-bind the chart to the (now swapped) source columns ---
$co = $ws.ChartObjects().Item(1)
$co.Chart.SetSourceData($ws.Range("B2:C10"))

# --- Widen column B to fit the new, longer "Measured Voltage [mV]" header ---
$ws.Columns.Item(2).ColumnWidth = 16.819820

# --- Add a new (currently blank) column D, matching column B's formatting ---
$ws.Range("B1:B10").Copy()
$ws.Range("D1:D10").PasteSpecial(-4122)
$ws.Columns.Item(4).ColumnWidth = 16.819820

# --- Move the selection cursor ---
$ws.Range("F8").Select()
